$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header labels for the "process status period" columns.
# Old layout (AX1:BF1): ВП, ПО, АЦ, СЦ, Кк, КД, НПР, Р, Зя
# New layout (AX1:BO1): Вб, SВб, ПО, SПО, АЦ, SАЦ, СЦ, SСЦ, ПцКк, SПцКк, Кк, ПцКД, SКД, ПцР, SПцР, Р, P-, Зя
$headers = @(
    "Вб", "SВб", "ПО", "SПО", "АЦ", "SАЦ", "СЦ", "SСЦ", "ПцКк",
    "SПцКк", "Кк", "ПцКД", "SКД", "ПцР", "SПцР", "Р", "P-", "Зя"
)

# Extend the existing header formatting (style) across the new, wider range
# before writing values, so every new header cell keeps the same look.
$ws.Range("AX1").Copy()
$ws.Range("AX1:BO1").PasteSpecial(-4122)

$startCol = 50 # column AX
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $startCol + $i).Value = $headers[$i]
}

# Update dimension/selection bookkeeping to point at the new last cell.
[void]$ws.Range("BO2").Select()
